$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the style of the other header cells (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill H2:H5 with 0
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0

# Fill H6 with 1
$ws.Range("H6").Value = 1
